$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4713456666666667
$ws.Range("H2").Value = 1.414037
$ws.Range("I2").Value = 0.003417883648355422
$ws.Range("J2").Value = 0.003417883648355422
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4059903333333333
$ws.Range("N2").Value = 1.217971
$ws.Range("O2").Value = 0.1286719387498537
$ws.Range("P2").Value = 0.1286719387498537
$ws.Range("Q2").Value = 0.1913617843252222
$ws.Range("R2").Value = 1.722256058927
$ws.Range("S2").Value = 0.0004397857154553154
$ws.Range("T2").Value = 0.0004397857154553154

# Row 3
$ws.Range("G3").Value = 0.4713456666666667
$ws.Range("H3").Value = 1.414037
$ws.Range("I3").Value = 0.003417883648355422
$ws.Range("J3").Value = 0.003417883648355422
$ws.Range("O3").Value = 0.4438927336444353
$ws.Range("P3").Value = 0.4438927336444353
$ws.Range("Q3").Value = 0.6601602990092222
$ws.Range("R3").Value = 5.941442691083
$ws.Range("S3").Value = 0.001517173715947104
$ws.Range("T3").Value = 0.001517173715947104

# Row 4
$ws.Range("G4").Value = 0.4713456666666667
$ws.Range("H4").Value = 1.414037
$ws.Range("I4").Value = 0.003417883648355422
$ws.Range("J4").Value = 0.003417883648355422
$ws.Range("M4").Value = 1.348659333333333
$ws.Range("N4").Value = 4.045978
$ws.Range("O4").Value = 0.4274353276057111
$ws.Range("P4").Value = 0.4274353276057111
$ws.Range("Q4").Value = 0.6356847325762222
$ws.Range("R4").Value = 5.721162593186
$ws.Range("S4").Value = 0.001460924216953003
$ws.Range("T4").Value = 0.001460924216953003

# Row 5
$ws.Range("I5").Value = 0.6765273260648147
$ws.Range("J5").Value = 0.6765273260648146
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4059903333333333
$ws.Range("N5").Value = 1.217971
$ws.Range("O5").Value = 0.1286719387498537
$ws.Range("P5").Value = 0.1286719387498537
$ws.Range("Q5").Value = 37.87767214452344
$ws.Range("R5").Value = 340.899049300711
$ws.Range("S5").Value = 0.08705008266201414
$ws.Range("T5").Value = 0.08705008266201412

# Row 6
$ws.Range("I6").Value = 0.6765273260648147
$ws.Range("J6").Value = 0.6765273260648146
$ws.Range("O6").Value = 0.4438927336444353
$ws.Range("P6").Value = 0.4438927336444353
$ws.Range("S6").Value = 0.3003055641520708
$ws.Range("T6").Value = 0.3003055641520708

# Row 7
$ws.Range("I7").Value = 0.6765273260648147
$ws.Range("J7").Value = 0.6765273260648146
$ws.Range("M7").Value = 1.348659333333333
$ws.Range("N7").Value = 4.045978
$ws.Range("O7").Value = 0.4274353276057111
$ws.Range("P7").Value = 0.4274353276057111
$ws.Range("Q7").Value = 125.8258432983664
$ws.Range("R7").Value = 1132.432589685298
$ws.Range("S7").Value = 0.2891716792507298
$ws.Range("T7").Value = 0.2891716792507297

# Row 8
$ws.Range("G8").Value = 44.13738266666667
$ws.Range("H8").Value = 132.412148
$ws.Range("I8").Value = 0.3200547902868299
$ws.Range("J8").Value = 0.3200547902868299
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4059903333333333
$ws.Range("N8").Value = 1.217971
$ws.Range("O8").Value = 0.1286719387498537
$ws.Range("P8").Value = 0.1286719387498537
$ws.Range("Q8").Value = 17.91935070130089
$ws.Range("R8").Value = 161.274156311708
$ws.Range("S8").Value = 0.04118207037238425
$ws.Range("T8").Value = 0.04118207037238425

# Row 9
$ws.Range("G9").Value = 44.13738266666667
$ws.Range("H9").Value = 132.412148
$ws.Range("I9").Value = 0.3200547902868299
$ws.Range("J9").Value = 0.3200547902868299
$ws.Range("O9").Value = 0.4438927336444353
$ws.Range("P9").Value = 0.4438927336444353
$ws.Range("Q9").Value = 61.81821495203689
$ws.Range("R9").Value = 556.363934568332
$ws.Range("S9").Value = 0.1420699957764174
$ws.Range("T9").Value = 0.1420699957764174

# Row 10
$ws.Range("G10").Value = 44.13738266666667
$ws.Range("H10").Value = 132.412148
$ws.Range("I10").Value = 0.3200547902868299
$ws.Range("J10").Value = 0.3200547902868299
$ws.Range("M10").Value = 1.348659333333333
$ws.Range("N10").Value = 4.045978
$ws.Range("O10").Value = 0.4274353276057111
$ws.Range("P10").Value = 0.4274353276057111
$ws.Range("Q10").Value = 59.52629308230488
$ws.Range("R10").Value = 535.736637740744
$ws.Range("S10").Value = 0.1368027241380283
$ws.Range("T10").Value = 0.1368027241380283
